$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for Goldschmidt (row 44), Labarga (row 58), White (row 113).
# Delete from bottom to top so earlier row numbers remain valid.
$ws.Rows.Item(113).Delete() | Out-Null
$ws.Rows.Item(58).Delete() | Out-Null
$ws.Rows.Item(44).Delete() | Out-Null
